$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("B2")
$c.Value = "''Bacteroides_cellulosilyticus_DSM_14838.mat'"
$c.Style = "Normal"
$ws.Range("C2").Value = 0.016

$c = $ws.Range("B3")
$c.Value = "''Bacteroides_coprocola_M16_DSM_17136.mat'"
$c.Style = "Normal"
$ws.Range("C3").Value = 0.003

$c = $ws.Range("B4")
$c.Value = "''Bacteroides_fluxus_YIT_12057.mat'"
$c.Style = "Normal"
$ws.Range("C4").Value = 0.048

$c = $ws.Range("B5")
$c.Value = "''Bacteroides_oleiciplenus_YIT_12058.mat'"
$c.Style = "Normal"
$ws.Range("C5").Value = 0.04

$c = $ws.Range("B6")
$c.Value = "''Bacteroides_ovatus_ATCC_8483.mat'"
$c.Style = "Normal"
$ws.Range("C6").Value = 0.23

$c = $ws.Range("B7")
$c.Value = "''Bacteroides_salyersiae_WAL_10018.mat'"
$c.Style = "Normal"
$ws.Range("C7").Value = 0.291

$c = $ws.Range("B8")
$c.Value = "''Bacteroides_stercoris_ATCC_43183.mat'"
$c.Style = "Normal"
$ws.Range("C8").Value = 0

$c = $ws.Range("B9")
$c.Value = "''Bacteroides_thetaiotaomicron_VPI_5482.mat'"
$c.Style = "Normal"
$ws.Range("C9").Value = 0

$c = $ws.Range("B10")
$c.Value = "''Bacteroides_uniformis_ATCC_8492.mat'"
$c.Style = "Normal"
$ws.Range("C10").Value = 0

$c = $ws.Range("B11")
$c.Value = "''Bacteroides_vulgatus_ATCC_8482.mat'"
$c.Style = "Normal"
$ws.Range("C11").Value = 0.372

$c = $ws.Range("B12")
$c.Value = "''Bifidobacterium_animalis_lactis_AD011.mat'"
$c.Style = "Normal"
$ws.Range("C12").Value = -0

$c = $ws.Range("B13")
$c.Value = "''Enterococcus_faecalis_OG1RF_ATCC_47077.mat'"
$c.Style = "Normal"
$ws.Range("C13").Value = 0

$c = $ws.Range("B14")
$c.Value = "''Flavonifractor_plautii_ATCC_29863.mat'"
$c.Style = "Normal"
$ws.Range("C14").Value = -0

$c = $ws.Range("B15")
$c.Value = "''Gordonibacter_pamelaeae_7_10_1_bT_DSM_19378.mat'"
$c.Style = "Normal"
$ws.Range("C15").Value = 0

$c = $ws.Range("B16")
$c.Value = "''Lactobacillus_plantarum_JDM1.mat'"
$c.Style = "Normal"
$ws.Range("C16").Value = -0

$c = $ws.Range("B17")
$c.Value = "''Odoribacter_laneus_YIT_12061.mat'"
$c.Style = "Normal"
$ws.Range("C17").Value = 0

$c = $ws.Range("B18")
$c.Value = "''Parabacteroides_distasonis_ATCC_8503.mat'"
$c.Style = "Normal"
$ws.Range("C18").Value = 0

$c = $ws.Range("B19")
$c.Value = "''Parabacteroides_johnsonii_DSM_18315.mat'"
$c.Style = "Normal"
$ws.Range("C19").Value = -0
